# Sync attendance_reports: reorder the "Recorded By" (column G) names so that
# the actual human reviewer's email (dnasr281@gmail.com / backup@backdoor.com)
# is listed before "System", instead of after it.
#
# Rule (derived from the diff): for every row whose G cell contains a comma
# separated list, if the list contains "dnasr281@gmail.com" or
# "backup@backdoor.com" then swap the first two comma-separated entries,
# leaving any further entries (e.g. a trailing lowercase duplicate "system")
# untouched. Rows that don't contain either of those addresses (e.g. the
# "System, admin@admin.com" only rows) are left unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Row + $used.Rows.Count - 1

$col = 7  # column G = "Recorded By"

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $val = $cell.Text

    if ([string]::IsNullOrEmpty($val)) {
        continue
    }
    if ($val.IndexOf(",") -lt 0) {
        continue
    }

    $needsSwap = ($val.Contains("dnasr281@gmail.com")) -or ($val.Contains("backup@backdoor.com"))
    if (-not $needsSwap) {
        continue
    }

    $rawParts = $val.Split(",")
    if ($rawParts.Count -lt 2) {
        continue
    }

    $parts = @()
    foreach ($p in $rawParts) {
        $parts += $p.Trim()
    }

    $tmp = $parts[0]
    $parts[0] = $parts[1]
    $parts[1] = $tmp

    $newVal = [string]::Join(", ", $parts)
    $cell.Value = $newVal
}
